$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Update Only" column (M) - header + "No" for each data row
$ws.Range("M1").Value = "Update Only"
$ws.Range("M2").Value = "No"
$ws.Range("M3").Value = "No"
$ws.Range("M4").Value = "No"
$ws.Range("M5").Value = "No"
$ws.Range("M6").Value = "No"
$ws.Range("M7").Value = "No"

# Match the formatting used by the rest of the imported offer rows (Normal 2
# cell style - Arial 11, no fill) for the new column.
$ws.Range("A2").Copy()
$ws.Range("M1:M7").PasteSpecial(-4122)

# The investor rows (6 & 7) that were previously styled slightly differently
# (A, B, I, J, L) are normalized to match the rest of the data rows.
$ws.Range("A6:B6").PasteSpecial(-4122)
$ws.Range("I6:J6").PasteSpecial(-4122)
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("I7:J7").PasteSpecial(-4122)
$ws.Range("L7").PasteSpecial(-4122)

# Restore selection/view state to the newly added column.
$ws.Activate() | Out-Null
$ws.Range("M3:M7").Select() | Out-Null
